# PANDUAN IKPK_D2.docx edit script
#
# 1. "Nyanyia Damai" (typo) -> "Nyanyian Damai", split across two runs
#    ("Nyanyian" / " Damai") matching the target OOXML.
# 2. "Bukusampul isi 38 lembar (menyesuaikan warna kelompok ) : NAMA, ..."
#    -> "Buku  isi 38 lembar: NAMA, ..." (drop the parenthetical remark,
#    split "Bukusampul" into "Buku" + " ").
# 3. "Sampul :" gains a new trailing run " Kertas Marmer (Menyesuaikan
#    Warna Kelompok)".
# 4. "Panduan Lengkap" (closing credit line) gains a new leading run
#    "(revisi 2) ".

$d = $word.ActiveDocument

function Split-RunAt($doc, $pos, $len) {
    # Force a hard run boundary at [$pos, $pos+$len) without altering the
    # visible formatting: flip Bold on then immediately back off. Each
    # assignment is a distinct edit, so the engine keeps the run split
    # even though the final Bold value matches the original.
    $r = $doc.Range($pos, $pos + $len)
    $orig = $r.Bold
    $r.Bold = 1
    $r2 = $doc.Range($pos, $pos + $len)
    $r2.Bold = $orig
}

# ---------------------------------------------------------------------
# Change 1: "Nyanyia Damai" -> "Nyanyian Damai" (split into 2 runs)
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Nyanyia Damai", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng.Find.Found) {
    $start = $rng.Start
    $rng.Text = "Nyanyian Damai"
    $firstLen = "Nyanyian".Length
    Split-RunAt $d $start $firstLen
}

# ---------------------------------------------------------------------
# Change 2: "Bukusampul" -> "Buku" + " "
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Bukusampul", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng.Find.Found) {
    $bukuStart = $rng.Start
    $bukuEnd = $rng.Start + 4
    $sampulEnd = $rng.End
    $sampulRng = $d.Range($bukuEnd, $sampulEnd)
    $sampulRng.Text = " "
    $bukuLen = 4
    Split-RunAt $d $bukuStart $bukuLen
}

# Change 2 (cont.): drop the "(menyesuaikan warna kelompok )" remark, so
# " isi 38 lembar (menyesuaikan warna kelompok ) : NAMA..." collapses to
# " isi 38 lembar: NAMA..."
$rng = $d.Content
$rng.Find.Execute(" (menyesuaikan warna kelompok ) ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng.Find.Found) {
    $rng.Text = ""
}

# ---------------------------------------------------------------------
# Change 3: "Sampul :" gains a new trailing run.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Sampul :", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng.Find.Found) {
    $insStart = $rng.End
    $addition = " Kertas Marmer (Menyesuaikan Warna Kelompok)"
    $addLen = $addition.Length
    $insPoint = $d.Range($insStart, $insStart)
    $insPoint.InsertAfter($addition)
    Split-RunAt $d $insStart $addLen
}

# ---------------------------------------------------------------------
# Change 4: "Panduan Lengkap" gains a new leading run "(revisi 2) ".
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Panduan Lengkap", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng.Find.Found) {
    $start = $rng.Start
    $prefix = "(revisi 2) "
    $prefixLen = $prefix.Length
    $rng.Text = $prefix + "Panduan Lengkap"
    Split-RunAt $d $start $prefixLen
}

Write-Output "edit complete"
